{"js": "// The diff removes every paragraph that follows the paragraph ending in\n// \"Il bianco, infine, ... i campi di testo.\" (i.e. the \"Prime fasi di test\n// dell'interfaccia\" heading and all subsequent paragraphs through \"In Figma\n// ... nell'applicazione finale.\"), replacing that whole block with a single\n// empty paragraph, right before the section properties.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph: the content to delete starts right after it.\nconst anchorText =\n  \"Il bianco, infine, viene sfruttato per evidenziare campi con i quali \u00e8 possibile interagire\";\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorText) !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\nif (anchorIndex === -1) {\n  throw new Error(\"Could not locate the anchor paragraph ('Il bianco, infine, ...').\");\n}\nif (anchorIndex + 1 >= paragraphs.items.length) {\n  throw new Error(\"No paragraphs found after the anchor paragraph to remove.\");\n}\n\n// OOXML for a completely bare, empty paragraph (no pPr/style/run at all).\nconst emptyParagraphOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body><w:p/></w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\n// Replace the first paragraph after the anchor (the \"Prime fasi di test\n// dell'interfaccia\" heading) with a bare empty paragraph - this keeps its\n// position as the direct successor of the anchor paragraph.\nconst firstParagraphToRemove = paragraphs.items[anchorIndex + 1];\nfirstParagraphToRemove.insertOoxml(emptyParagraphOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// Delete every remaining paragraph that used to follow (from the end of the\n// document backwards, to keep indices valid while deleting).\nparagraphs.load(\"items/text\");\nawait context.sync();\nfor (let i = paragraphs.items.length - 1; i >= anchorIndex + 2; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# The diff removes every paragraph that follows the paragraph ending in\n# \"Il bianco, infine, ... i campi di testo.\" (i.e. the \"Prime fasi di test\n# dell'interfaccia\" heading and all subsequent paragraphs through \"In Figma\n# ... nell'applicazione finale.\"), replacing that whole block with a single\n# empty paragraph, right before the section properties.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph: the content to delete starts right after it.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Il bianco, infine, viene sfruttato per evidenziare campi*\") {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not locate the anchor paragraph ('Il bianco, infine, ...').\"\n}\n\n$firstToRemove = $anchorIndex + 1\n$lastToRemove = $d.Paragraphs.Count\nif ($firstToRemove -gt $lastToRemove) {\n    throw \"No paragraphs found after the anchor paragraph to remove.\"\n}\n\n# Delete everything from the start of the paragraph right after the anchor\n# (keep that one paragraph itself, so one paragraph mark survives) through\n# the start of the last paragraph to remove's following paragraph, i.e.\n# delete all paragraphs strictly after the first one to remove.\nif ($lastToRemove -gt $firstToRemove) {\n    $deleteStart = $d.Paragraphs.Item($firstToRemove + 1).Range.Start\n    $deleteEnd = $d.Paragraphs.Item($lastToRemove).Range.End\n    $d.Range($deleteStart, $deleteEnd).Delete()\n}\n\n# Now clear the remaining \"first paragraph to remove\" down to a completely\n# bare paragraph: no text, no run, no paragraph formatting/style override.\n$para = $d.Paragraphs.Item($firstToRemove)\n$r = $para.Range\n[void]$r.MoveEnd(1, -1)\n$r.Delete()\n$para.Style = \"Normal\"\n"}
